$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Insert a new "2022-Q4" sheet, positioned right after "2022-Q3"
#    and before it in tab order, by duplicating "2022-Q3" (so all
#    styling/column widths/header formatting matches exactly) and
#    then overwriting its data with the Q4 numbers.
# ---------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item($q3.Index - 1)
$q4.Name = "2022-Q4"

# The copied sheet has 57 data rows (2022-Q3 had 57 funds); the
# 2022-Q4 snapshot only has 43, so drop the extra rows (45-58).
$q4.Range("A45:H58").EntireRow.Delete()

# Column A (rank index 0..42) is already correct after the copy -
# only B..H need new values. Force B:G to stay text (they hold
# numbers formatted as text in the source data), write everything
# in one shot via a 2-D array, then strip the temporary text format
# so the cells end up with no explicit style (matching the rest of
# the data rows).
$textRange = $q4.Range("B2:G44")
$textRange.NumberFormat = "@"

$arr = New-Object 'object[,]' 43,7
$arr[0,0]='012079'; $arr[0,1]='信澳新能源精选混合'; $arr[0,2]='40.55'; $arr[0,3]='93.85'; $arr[0,4]='9.55'; $arr[0,5]='3.8725'; $arr[0,6]=4
$arr[1,0]='006736'; $arr[1,1]='国投瑞银先进制造混合'; $arr[1,2]='38.57'; $arr[1,3]='94.27'; $arr[1,4]='8.18'; $arr[1,5]='3.1550'; $arr[1,6]=1
$arr[2,0]='007689'; $arr[2,1]='国投瑞银新能源混合A'; $arr[2,2]='37.33'; $arr[2,3]='94.76'; $arr[2,4]='8.17'; $arr[2,5]='3.0499'; $arr[2,6]=2
$arr[3,0]='007690'; $arr[3,1]='国投瑞银新能源混合C'; $arr[3,2]='30.69'; $arr[3,3]='94.76'; $arr[3,4]='8.17'; $arr[3,5]='2.5074'; $arr[3,6]=2
$arr[4,0]='012148'; $arr[4,1]='国投瑞银产业趋势混合A'; $arr[4,2]='25.47'; $arr[4,3]='93.75'; $arr[4,4]='8.19'; $arr[4,5]='2.0860'; $arr[4,6]=1
$arr[5,0]='012149'; $arr[5,1]='国投瑞银产业趋势混合C'; $arr[5,2]='15.53'; $arr[5,3]='93.75'; $arr[5,4]='8.19'; $arr[5,5]='1.2719'; $arr[5,6]=1
$arr[6,0]='012223'; $arr[6,1]='信澳成长精选混合A'; $arr[6,2]='10.50'; $arr[6,3]='87.86'; $arr[6,4]='4.98'; $arr[6,5]='0.5229'; $arr[6,6]=9
$arr[7,0]='512100'; $arr[7,1]='南方中证1000ETF'; $arr[7,2]='106.09'; $arr[7,3]='98.15'; $arr[7,4]='0.47'; $arr[7,5]='0.4986'; $arr[7,6]=1
$arr[8,0]='560010'; $arr[8,1]='广发中证1000ETF'; $arr[8,2]='67.21'; $arr[8,3]='98.32'; $arr[8,4]='0.48'; $arr[8,5]='0.3226'; $arr[8,6]=1
$arr[9,0]='159845'; $arr[9,1]='华夏中证1000ETF'; $arr[9,2]='62.11'; $arr[9,3]='98.70'; $arr[9,4]='0.48'; $arr[9,5]='0.2981'; $arr[9,6]=1
$arr[10,0]='005119'; $arr[10,1]='银华智荟内在价值灵活配置混合A'; $arr[10,2]='6.43'; $arr[10,3]='94.86'; $arr[10,4]='4.45'; $arr[10,5]='0.2861'; $arr[10,6]=7
$arr[11,0]='159629'; $arr[11,1]='富国中证1000ETF'; $arr[11,2]='59.33'; $arr[11,3]='99.34'; $arr[11,4]='0.48'; $arr[11,5]='0.2848'; $arr[11,6]=1
$arr[12,0]='159633'; $arr[12,1]='易方达中证1000ETF'; $arr[12,2]='58.64'; $arr[12,3]='98.77'; $arr[12,4]='0.48'; $arr[12,5]='0.2815'; $arr[12,6]=1
$arr[13,0]='610004'; $arr[13,1]='信澳中小盘混合'; $arr[13,2]='4.60'; $arr[13,3]='91.92'; $arr[13,4]='6.06'; $arr[13,5]='0.2788'; $arr[13,6]=5
$arr[14,0]='012224'; $arr[14,1]='信澳成长精选混合C'; $arr[14,2]='5.34'; $arr[14,3]='87.86'; $arr[14,4]='4.98'; $arr[14,5]='0.2659'; $arr[14,6]=9
$arr[15,0]='009859'; $arr[15,1]='银华乐享混合A'; $arr[15,2]='4.40'; $arr[15,3]='94.60'; $arr[15,4]='4.96'; $arr[15,5]='0.2182'; $arr[15,6]=9
$arr[16,0]='562800'; $arr[16,1]='嘉实中证稀有金属主题ETF'; $arr[16,2]='6.62'; $arr[16,3]='99.41'; $arr[16,4]='3.18'; $arr[16,5]='0.2105'; $arr[16,6]=8
$arr[17,0]='610006'; $arr[17,1]='信澳产业升级混合'; $arr[17,2]='4.01'; $arr[17,3]='90.01'; $arr[17,4]='4.34'; $arr[17,5]='0.1740'; $arr[17,6]=10
$arr[18,0]='159667'; $arr[18,1]='国泰中证机床ETF'; $arr[18,2]='3.49'; $arr[18,3]='99.21'; $arr[18,4]='4.71'; $arr[18,5]='0.1644'; $arr[18,6]=5
$arr[19,0]='290014'; $arr[19,1]='泰信现代服务业混合'; $arr[19,2]='2.40'; $arr[19,3]='77.48'; $arr[19,4]='5.20'; $arr[19,5]='0.1248'; $arr[19,6]=4
$arr[20,0]='013495'; $arr[20,1]='信澳产业优选一年持有混合A'; $arr[20,2]='1.92'; $arr[20,3]='79.53'; $arr[20,4]='4.15'; $arr[20,5]='0.0797'; $arr[20,6]=10
$arr[21,0]='011939'; $arr[21,1]='博时新能源汽车主题混合C'; $arr[21,2]='1.56'; $arr[21,3]='88.58'; $arr[21,4]='4.44'; $arr[21,5]='0.0693'; $arr[21,6]=6
$arr[22,0]='290008'; $arr[22,1]='泰信发展主题混合'; $arr[22,2]='1.18'; $arr[22,3]='84.83'; $arr[22,4]='5.63'; $arr[22,5]='0.0664'; $arr[22,6]=4
$arr[23,0]='015687'; $arr[23,1]='银华乐享混合C'; $arr[23,2]='1.29'; $arr[23,3]='94.60'; $arr[23,4]='4.96'; $arr[23,5]='0.0640'; $arr[23,6]=9
$arr[24,0]='013104'; $arr[24,1]='博时新能源主题混合C'; $arr[24,2]='2.21'; $arr[24,3]='86.45'; $arr[24,4]='2.80'; $arr[24,5]='0.0619'; $arr[24,6]=8
$arr[25,0]='013103'; $arr[25,1]='博时新能源主题混合A'; $arr[25,2]='2.00'; $arr[25,3]='86.45'; $arr[25,4]='2.80'; $arr[25,5]='0.0560'; $arr[25,6]=8
$arr[26,0]='011938'; $arr[26,1]='博时新能源汽车主题混合A'; $arr[26,2]='1.22'; $arr[26,3]='88.58'; $arr[26,4]='4.44'; $arr[26,5]='0.0542'; $arr[26,6]=6
$arr[27,0]='013943'; $arr[27,1]='华宝中证稀有金属指数增强C'; $arr[27,2]='1.22'; $arr[27,3]='94.08'; $arr[27,4]='4.24'; $arr[27,5]='0.0517'; $arr[27,6]=7
$arr[28,0]='159608'; $arr[28,1]='广发中证稀有金属ETF'; $arr[28,2]='1.59'; $arr[28,3]='98.42'; $arr[28,4]='3.19'; $arr[28,5]='0.0507'; $arr[28,6]=8
$arr[29,0]='159663'; $arr[29,1]='华夏中证机床ETF'; $arr[29,2]='1.01'; $arr[29,3]='97.42'; $arr[29,4]='4.65'; $arr[29,5]='0.0470'; $arr[29,6]=6
$arr[30,0]='016262'; $arr[30,1]='银华智荟内在价值灵活配置混合C'; $arr[30,2]='0.82'; $arr[30,3]='94.86'; $arr[30,4]='4.45'; $arr[30,5]='0.0365'; $arr[30,6]=7
$arr[31,0]='561800'; $arr[31,1]='华富中证稀有金属主题ETF'; $arr[31,2]='0.94'; $arr[31,3]='99.79'; $arr[31,4]='3.20'; $arr[31,5]='0.0301'; $arr[31,6]=7
$arr[32,0]='560110'; $arr[32,1]='汇添富中证1000ETF'; $arr[32,2]='5.61'; $arr[32,3]='93.96'; $arr[32,4]='0.45'; $arr[32,5]='0.0252'; $arr[32,6]=1
$arr[33,0]='013942'; $arr[33,1]='华宝中证稀有金属指数增强A'; $arr[33,2]='0.35'; $arr[33,3]='94.08'; $arr[33,4]='4.24'; $arr[33,5]='0.0148'; $arr[33,6]=7
$arr[34,0]='001572'; $arr[34,1]='嘉合磐石混合C'; $arr[34,2]='0.39'; $arr[34,3]='37.90'; $arr[34,4]='2.94'; $arr[34,5]='0.0115'; $arr[34,6]=8
$arr[35,0]='003646'; $arr[35,1]='创金合信中证1000指数增强A'; $arr[35,2]='0.52'; $arr[35,3]='92.80'; $arr[35,4]='1.53'; $arr[35,5]='0.0080'; $arr[35,6]=1
$arr[36,0]='013496'; $arr[36,1]='信澳产业优选一年持有混合C'; $arr[36,2]='0.18'; $arr[36,3]='79.53'; $arr[36,4]='4.15'; $arr[36,5]='0.0075'; $arr[36,6]=10
$arr[37,0]='516300'; $arr[37,1]='华泰柏瑞中证1000ETF'; $arr[37,2]='1.57'; $arr[37,3]='98.94'; $arr[37,4]='0.48'; $arr[37,5]='0.0075'; $arr[37,6]=1
$arr[38,0]='003647'; $arr[38,1]='创金合信中证1000指数增强C'; $arr[38,2]='0.43'; $arr[38,3]='92.80'; $arr[38,4]='1.53'; $arr[38,5]='0.0066'; $arr[38,6]=1
$arr[39,0]='159918'; $arr[39,1]='嘉实中创400ETF'; $arr[39,2]='0.59'; $arr[39,3]='98.55'; $arr[39,4]='0.73'; $arr[39,5]='0.0043'; $arr[39,6]=4
$arr[40,0]='001571'; $arr[40,1]='嘉合磐石混合A'; $arr[40,2]='0.09'; $arr[40,3]='37.90'; $arr[40,4]='2.94'; $arr[40,5]='0.0026'; $arr[40,6]=8
$arr[41,0]='162413'; $arr[41,1]='华宝中证1000指数A'; $arr[41,2]='0.40'; $arr[41,3]='91.60'; $arr[41,4]='0.44'; $arr[41,5]='0.0018'; $arr[41,6]=1
$arr[42,0]='016033'; $arr[42,1]='华宝中证1000指数C'; $arr[42,2]='0.09'; $arr[42,3]='91.60'; $arr[42,4]='0.44'; $arr[42,5]='0.0004'; $arr[42,6]=1

$q4.Range("B2:H44").Value = $arr
$textRange.ClearFormats()

# ---------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for
#    2022-Q4 and shift the existing quarters down by one row.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A2").EntireRow.Insert()
$summary.Range("A2:D2").ClearFormats()

# Re-apply the same style used by the other index cells in column A
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 43
$summary.Range("D2").Value = 20.63

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 57
$summary.Range("D3").Value = 27.27

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 90
$summary.Range("D4").Value = 28.31

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 38
$summary.Range("D5").Value = 27.56

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 17
$summary.Range("D6").Value = 16.16

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 25
$summary.Range("D7").Value = 23.11
